$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so numeric-looking strings
# (e.g. "11.01", "1.00") are preserved exactly as text, matching the source
# workbook where these cells are inline/shared strings, not numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "67.307.63"
$ws.Range("E2").Value = "  -0.94%  "

# Row 3
$ws.Range("D3").Value = "2.608.10"
$ws.Range("E3").Value = "  -0.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "590.11"
$ws.Range("E5").Value = "  -0.94%  "

# Row 6
$ws.Range("D6").Value = "149.62"
$ws.Range("E6").Value = "  -3.71%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.53%  "

# Row 9
$ws.Range("D9").Value = "2.608.24"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -4.16%  "

# Row 11
$ws.Range("E11").Value = "  -0.14%  "

# Row 12
$ws.Range("E12").Value = "  -1.30%  "

# Row 13
$ws.Range("E13").Value = "  -3.02%  "

# Row 14
$ws.Range("D14").Value = "27.30"
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("D15").Value = "3.080.67"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("E16").Value = "  -4.74%  "

# Row 17
$ws.Range("D17").Value = "67.208.79"
$ws.Range("E17").Value = "  -0.87%  "

# Row 18
$ws.Range("D18").Value = "2.608.85"
$ws.Range("E18").Value = "  -0.54%  "

# Row 19
$ws.Range("D19").Value = "366.73"
$ws.Range("E19").Value = "  -0.35%  "

# Row 20
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").Value = "  -1.81%  "

# Row 21
$ws.Range("D21").Value = "7.34"
$ws.Range("E21").Value = "  -4.14%  "

# Row 22
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("D23").Value = "4.82"
$ws.Range("E23").Value = "  -2.53%  "

# Row 24
$ws.Range("E24").Value = "  -0.63%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "9.94"
$ws.Range("E26").Value = "  +0.87%  "

# Row 27
$ws.Range("D27").Value = "67.38"
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("D28").Value = "2.743.74"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.26%  "

# Row 30
$ws.Range("D30").Value = "579.38"
$ws.Range("E30").Value = "  +0.18%  "

# Row 31
$ws.Range("E31").Value = "  -5.55%  "

# Row 32
$ws.Range("E32").Value = "  -5.01%  "

# Row 33
$ws.Range("D33").Value = "7.63"
$ws.Range("E33").Value = "  -3.89%  "

# Row 34
$ws.Range("E34").Value = "  -3.26%  "

# Row 35
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.07%  "

# Row 36
$ws.Range("D36").Value = "0.124"
$ws.Range("E36").Value = "  -6.39%  "

# Row 37
$ws.Range("D37").Value = "1.49"
$ws.Range("E37").Value = "  -2.71%  "

# Row 38
$ws.Range("D38").Value = "155.81"
$ws.Range("E38").Value = "  -2.04%  "

# Row 39
$ws.Range("D39").Value = "18.96"
$ws.Range("E39").Value = "  -2.19%  "

# Row 40
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  +0.95%  "

# Row 41
$ws.Range("E41").Value = "  -1.69%  "

# Row 42
$ws.Range("E42").Value = "  -2.71%  "

# Row 43
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "16.79"
$ws.Range("E43").Value = "  +2.24%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.49"
$ws.Range("E44").Value = "  -3.18%  "

# Row 45
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.15%  "

# Row 46
$ws.Range("D46").Value = "154.43"
$ws.Range("E46").Value = "  -1.30%  "

# Row 47
$ws.Range("D47").Value = "0.0₆0289"
$ws.Range("E47").Value = "  +0.39%  "

# Row 48
$ws.Range("E48").Value = "  -1.15%  "

# Row 49
$ws.Range("D49").Value = "0.0784"
$ws.Range("E49").Value = "  +0.09%  "

# Row 50
$ws.Range("E50").Value = "  -2.38%  "

# Row 51
$ws.Range("D51").Value = "21.32"
$ws.Range("E51").Value = "  +2.19%  "

# Reset style on the Price column back to Normal so no stray
# number-format/style index is left applied to these cells.
$priceRange.Style = "Normal"
